# Insert a new statistics row for "BMI (Kg/m2)" right above the existing
# "Lingkar Perut (cm)" row (i.e. as the new row 5), pushing the remaining
# rows (Lingkar Perut, Lingkar Leher, Terbangun x2, Durasi tidur) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 5; everything currently in row 5 onward shifts down.
$ws.Rows.Item(5).Insert()

# Copy the formatting from the (now shifted) header cell in A6 onto the new
# A5 cell so it keeps the bold/bordered/centered label style used by the
# other row headers in column A.
$ws.Range("A6").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# Fill in the new BMI row's label and statistics.
$ws.Range("A5").Value = "BMI (Kg/m2)"
$ws.Range("B5").Value = 9.036860879904875
$ws.Range("C5").Value = 178.8139429606157
$ws.Range("D5").Value = 29.40742939196596
$ws.Range("E5").Value = 10.18163704002391
